$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Absent" column (H) values that were left blank/zero while
# forming the consolidated report.
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 0
